{"js": "// Apply the three-digit x one-digit multiplication problem updates.\n// Each entry is [oldExpression, newExpression]; all values are unique in the\n// document so a direct search+replace per pair is unambiguous.\nconst replacements = [\n  [\"149\u00d77=\", \"531\u00d79=\"],\n  [\"731\u00d74=\", \"864\u00d76=\"],\n  [\"670\u00d79=\", \"415\u00d76=\"],\n  [\"659\u00d79=\", \"938\u00d72=\"],\n  [\"982\u00d75=\", \"590\u00d75=\"],\n  [\"760\u00d73=\", \"996\u00d77=\"],\n  [\"840\u00d78=\", \"860\u00d78=\"],\n  [\"162\u00d73=\", \"161\u00d78=\"],\n  [\"556\u00d76=\", \"599\u00d77=\"],\n  [\"300\u00d76=\", \"692\u00d74=\"],\n  [\"562\u00d79=\", \"164\u00d73=\"],\n  [\"603\u00d77=\", \"237\u00d74=\"],\n  [\"829\u00d73=\", \"838\u00d74=\"],\n  [\"148\u00d78=\", \"536\u00d74=\"],\n  [\"783\u00d78=\", \"923\u00d77=\"],\n  [\"760\u00d78=\", \"837\u00d72=\"],\n  [\"477\u00d77=\", \"914\u00d72=\"],\n  [\"455\u00d76=\", \"985\u00d77=\"],\n  [\"192\u00d77=\", \"444\u00d72=\"],\n  [\"208\u00d76=\", \"264\u00d75=\"],\n  [\"471\u00d76=\", \"501\u00d72=\"],\n  [\"574\u00d78=\", \"861\u00d72=\"],\n  [\"508\u00d77=\", \"924\u00d75=\"],\n  [\"880\u00d76=\", \"624\u00d79=\"],\n  [\"108\u00d76=\", \"509\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the three-digit x one-digit multiplication problem updates.\n# Each entry is old/new expression text; all values are unique in the\n# document so Find/Replace per pair (whole document, match case) is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('149\u00d77=', '531\u00d79='),\n    @('731\u00d74=', '864\u00d76='),\n    @('670\u00d79=', '415\u00d76='),\n    @('659\u00d79=', '938\u00d72='),\n    @('982\u00d75=', '590\u00d75='),\n    @('760\u00d73=', '996\u00d77='),\n    @('840\u00d78=', '860\u00d78='),\n    @('162\u00d73=', '161\u00d78='),\n    @('556\u00d76=', '599\u00d77='),\n    @('300\u00d76=', '692\u00d74='),\n    @('562\u00d79=', '164\u00d73='),\n    @('603\u00d77=', '237\u00d74='),\n    @('829\u00d73=', '838\u00d74='),\n    @('148\u00d78=', '536\u00d74='),\n    @('783\u00d78=', '923\u00d77='),\n    @('760\u00d78=', '837\u00d72='),\n    @('477\u00d77=', '914\u00d72='),\n    @('455\u00d76=', '985\u00d77='),\n    @('192\u00d77=', '444\u00d72='),\n    @('208\u00d76=', '264\u00d75='),\n    @('471\u00d76=', '501\u00d72='),\n    @('574\u00d78=', '861\u00d72='),\n    @('508\u00d77=', '924\u00d75='),\n    @('880\u00d76=', '624\u00d79='),\n    @('108\u00d76=', '509\u00d73='),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($null, $true, $null, $null, $null, $null, $true, $null, $null, $null, 2) | Out-Null\n}\n"}
